$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt Diagramm")

# --- Update team member names (E2, I2, M2, Q2) ---
$ws.Range("E2").Value = "Leon Kranner"
$ws.Range("I2").Value = "Marco Kuner"
$ws.Range("M2").Value = "David Vollmer"
$ws.Range("Q2").Value = "Marcel Wagner"

# --- Rename the "Zeiträume" header to "Kalenderwoche" (merged D3:O3) ---
$ws.Range("D3").Value = "Kalenderwoche"

# --- Unmerge the name cells in row 2, then restore left/vertical-center alignment ---
$ws.Range("E2:F2").UnMerge()
$ws.Range("I2:J2").UnMerge()
$ws.Range("M2:N2").UnMerge()
$ws.Range("Q2:R2").UnMerge()

$ws.Range("E2:F2").HorizontalAlignment = 1
$ws.Range("E2:F2").VerticalAlignment = -4108
$ws.Range("I2:J2").HorizontalAlignment = 1
$ws.Range("I2:J2").VerticalAlignment = -4108
$ws.Range("M2:N2").HorizontalAlignment = 1
$ws.Range("M2:N2").VerticalAlignment = -4108
$ws.Range("Q2:R2").HorizontalAlignment = 1
$ws.Range("Q2:R2").VerticalAlignment = -4108

# --- Shift the calendar week numbers in row 4 from 1-16 to 12-27 ---
for ($col = 4; $col -le 19; $col++) {
    $ws.Cells.Item(4, $col).Value = $col + 8
}

# --- Update the active selection to match the saved view ---
$ws.Range("T4").Select() | Out-Null
